# FeSources.xlsx edit script
# - Corrects an erroneous "breakdown" hyperlink/text that was left in the
#   Employer Skills Survey row of the Sources sheet (it mistakenly carried
#   text/link for the ONS Labour Market Profile instead of being empty).
# - Adds a new row describing the list of published Local Skills
#   Improvement Plans (LSIPs), with its linked "Employer representative
#   bodies" source.
# - Leaves the Sources sheet as the active tab/selection (was Tools).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Tools"
$ws2 = $wb.Worksheets.Item(2)   # "Sources"

# --- Fix the erroneous breakdown hyperlink/text on the ESS row (old C11) ---
$ws2.Hyperlinks.Delete()
$ws2.Range("C11").ClearContents()

# --- Insert a new row above the ESS row for the LSIP / ERB source ---
$ws2.Rows.Item(11).Insert()

# Copy formatting from the row above so the new row matches the sheet's
# existing look (border/alignment/wrap style), without creating a new style.
$ws2.Range("A10").Copy()
$ws2.Range("A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data column first isn't required, but Source (B) is entered before Data (A)
# so new shared strings land in the same order as the authored workbook.
$ws2.Range("B11").Value = "<a href='https://www.gov.uk/government/publications/designated-employer-representative-bodies/notice-of-designated-employer-representative-bodies'>Employer representative bodies</a>"
$ws2.Range("A11").Value = "List of links to published Local skills improvement plans. "
$ws2.Rows.Item(11).RowHeight = 29

# --- Keep the page setup explicit on the Sources sheet ---
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- Make "Sources" the active/selected tab, with F10 selected ---
$ws2.Activate()
$ws2.Range("F10").Select()
